$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "pander(table_forecasts)*") {
        $p.Range.Delete()
        break
    }
}
